$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("UI")
$ws2 = $wb.Worksheets.Item("Resistances")
$ws3 = $wb.Worksheets.Item("Feuil3")

# ---- Fill in Feuil3 (sheet3) data: values & formulas ----

$ws3.Range('A6').Value = 'réf X'
$ws3.Range('B6').Value = 'X0'
$ws3.Range('C6').Value = 0.0
$ws3.Range('D6').Value = 0.0
$ws3.Range('E6').Value = 0.10000000000000001
$ws3.Range('F6').Value = 0.10000000000000001
$ws3.Range('A7').Value = 'réf Y'
$ws3.Range('B7').Value = 'Y0'
$ws3.Range('C7').Value = -0.42999999999999999
$ws3.Range('D7').Value = -0.41999999999999998
$ws3.Range('E7').Value = 0.00231
$ws3.Range('F7').Value = 0.10000000000000001
$ws3.Range('A8').Value = 'mes X'
$ws3.Range('B8').Value = 'x1'
$ws3.Range('C8').Value = 450.0
$ws3.Range('D8').Value = 450.0
$ws3.Range('E8').Value = 23.5
$ws3.Range('F8').Value = 23.5
$ws3.Range('A9').Value = 'mes Y'
$ws3.Range('B9').Value = 'Y1'
$ws3.Range('C9').Value = 448.33999999999997499
$ws3.Range('D9').Value = 448.52999999999997272
$ws3.Range('E9').Value = 23.48510999999999882
$ws3.Range('F9').Value = 23.49780000000000157
$ws3.Range('A10').Value = 'Dréf'
$ws3.Range('B10').Value = 'Dx'
$ws3.Range('C10').Formula = '=C8-C6'
$ws3.Range('D10').Formula = '=D8-D6'
$ws3.Range('F10').Formula = '=F8-F6'
$ws3.Range('A11').Value = 'Dmes'
$ws3.Range('B11').Value = 'Dy'
$ws3.Range('C11').Formula = '=C9-C7'
$ws3.Range('D11').Formula = '=D9-D7'
$ws3.Range('F11').Formula = '=F9-F7'
$ws3.Range('A12').Value = 'offset'
$ws3.Range('B12').Value = 'b'''
$ws3.Range('C12').Formula = '=$C6-C7'
$ws3.Range('D12').Formula = '=$C6-D7'
$ws3.Range('E12').Formula = '=$D6-E7'
$ws3.Range('F12').Formula = '=$D6-F7'
$ws3.Range('A13').Value = 'pente mes'
$ws3.Range('B13').Value = 'a'''
$ws3.Range('C13').Formula = '=C11/C10'
$ws3.Range('D13').Formula = '=D11/D10'
$ws3.Range('E13').Formula = '=(E9-E7)/($D8-$D6)'
$ws3.Range('F13').Formula = '=F11/F10'
$ws3.Range('A14').Value = 'gain origine'
$ws3.Range('B14').Value = 'g'''
$ws3.Range('C14').Value = 1.0
$ws3.Range('D14').Value = 1.0
$ws3.Range('E14').Value = 1.00387500000000007
$ws3.Range('F14').Value = 1.00022999999999995
$ws3.Range('A15').Value = 'gain corrigé'
$ws3.Range('B15').Value = 'g''/a'''
$ws3.Range('C15').Formula = '=C14/C13'
$ws3.Range('D15').Formula = '=D14/D13'
$ws3.Range('E15').Formula = '=E14/E13'
$ws3.Range('F15').Formula = '=F14/F13'
$ws3.Range('B17').Value = 'z'
$ws3.Range('C17').Value = 0.0
$ws3.Range('E17').Value = 0.0098
$ws3.Range('F17').Value = 0.0098
$ws3.Range('B18').Value = 'z-b'''
$ws3.Range('C18').Formula = '=C17-C12'
$ws3.Range('E18').Formula = '=E17-E12'
$ws3.Range('F18').Formula = '=F17-F12'
$ws3.Range('A19').Value = 'correction Z'
$ws3.Range('C19').Formula = '=C18*1/-486'
$ws3.Range('C21').Formula = '=C11*C15'
$ws3.Range('C22').Formula = '=0.447979999999999/10000'

# ---- Apply cell styles (reuse existing style entries via PasteSpecial Formats) ----
$ws1.Range('D15').Copy()
$ws3.Range('C16').PasteSpecial(-4122)
$ws3.Range('E16').PasteSpecial(-4122)
$ws3.Range('F16').PasteSpecial(-4122)
$ws1.Range('D5').Copy()
$ws3.Range('C6').PasteSpecial(-4122)
$ws3.Range('D6').PasteSpecial(-4122)
$ws3.Range('E6').PasteSpecial(-4122)
$ws3.Range('C8').PasteSpecial(-4122)
$ws3.Range('D8').PasteSpecial(-4122)
$ws3.Range('E8').PasteSpecial(-4122)
$ws1.Range('G3').Copy()
$ws3.Range('A6').PasteSpecial(-4122)
$ws3.Range('A7').PasteSpecial(-4122)
$ws3.Range('F7').PasteSpecial(-4122)
$ws3.Range('A8').PasteSpecial(-4122)
$ws3.Range('A9').PasteSpecial(-4122)
$ws3.Range('F9').PasteSpecial(-4122)
$ws3.Range('A10').PasteSpecial(-4122)
$ws3.Range('A11').PasteSpecial(-4122)
$ws3.Range('A12').PasteSpecial(-4122)
$ws3.Range('A13').PasteSpecial(-4122)
$ws3.Range('A14').PasteSpecial(-4122)
$ws3.Range('A15').PasteSpecial(-4122)
$ws3.Range('A19').PasteSpecial(-4122)
$ws3.Range('F19').PasteSpecial(-4122)
$ws1.Range('G5').Copy()
$ws3.Range('F6').PasteSpecial(-4122)
$ws3.Range('F8').PasteSpecial(-4122)

# ---- Column width for column A on Feuil3 ----
$ws3.Columns.Item(1).ColumnWidth = 14.7

# ---- View / selection changes ----
$ws1.Range("C5:G18").Select()
$ws2.Range("B13").Select()
$ws3.Range("C21").Select()

